# Events.xlsx edit: insert new "e005a" event row before the existing "e006" row
# (single-day-of-battle intro event), update row heights, and refresh the
# worksheet selection/scroll position to match the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7; this shifts every row from the old row 7
# (e006 "Combat Calendar Check") down by one, and the new row inherits the
# formatting (styles) of the row above it, same as interactive Excel.
$ws.Rows("7:7").Insert()

# Populate the new row with the "e005a" event key and its body text.
$ws.Range("A7").Value = "e005a"
$ws.Range("B7").Value = "<Bold>e005a Single Day of Battle</Bold> `n<LineBreak/><LineBreak/>`nYou elected to only fight one day of battle. The game is won or lost based on victory points at end of day per `n<InlineUIContainer><Button Content='r6.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. `n<LineBreak/><LineBreak/>`nClick buttons to get day and tank you want. The tank choice is limited by date shown on `n<InlineUIContainer><Button Content='Replacement' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. Click image when ready to begin.`n<LineBreak/><LineBreak/>"

# Match the authored row height for the new row (120pt, vs. the old 210pt
# that row 7 had before the insert -- that content is now on row 8).
$ws.Rows("7:7").RowHeight = 120

# Update the view: scroll so row 4 is the top visible row, and select B7 like
# the author left it selected after adding the new content.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select()
